$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New check-in data to concat onto the existing sheet, starting at row 2
# (the prior single data row is overwritten as part of the same paste).
# All values are plain text, matching the sheet's existing inline-string
# formatting (Student ID is a text code, not a numeric value).
$data = @(
    @("2151034", "Juanna", "2021-09-30 12:00:00"),
    @("2151034", "Juanna", "2021-09-30 12:00:00"),
    @("2151034", "Juanna", "2021-09-30 12:00:00"),
    @("2151034", "Juanna", "2021-09-30 12:00:00"),
    @("2151034", "Juanna", "2021-09-30 12:00:00"),
    @("2151034", "Juanna", "15/04/2024 10:46:41"),
    @("2151034", "Juanna", "15/04/2024 10:47:13")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i

    # Student ID looks numeric ("2151034") but must stay text, like the
    # original "0003694140" value did - force text formatting before
    # assigning, then drop back to the Normal style so no stray number
    # format is left behind on the cell.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
